$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2731
$ws.Range("E2").Value = 425
$ws.Range("F2").Value = 425
$ws.Range("G2").Value = 526
$ws.Range("H2").Value = 449
$ws.Range("I2").Value = 443
$ws.Range("J2").Value = 6
$ws.Range("K2").Value = 6770
$ws.Range("L2").Value = 1198
$ws.Range("M2").Value = 5573
$ws.Range("N2").Value = 5470
$ws.Range("O2").Value = 102
$ws.Range("P2").Value = 452
$ws.Range("Q2").Value = 126
$ws.Range("R2").Value = -48
$ws.Range("S2").Value = 25
$ws.Range("T2").Value = 138
$ws.Range("U2").Value = -12
$ws.Range("V2").Value = 406
$ws.Range("W2").Value = 15.56
$ws.Range("X2").Value = 16.45
$ws.Range("Y2").Value = 8.36
$ws.Range("Z2").Value = 6.87
$ws.Range("AA2").Value = 21.5
$ws.Range("AB2").Value = 1233.4
$ws.Range("AC2").Value = 5654
$ws.Range("AD2").Value = 7.58
$ws.Range("AE2").Value = 84266
$ws.Range("AF2").Value = 0.51
$ws.Range("AG2").Value = 1200
$ws.Range("AH2").Value = 2.8
$ws.Range("AI2").Value = 17.58
$ws.Range("AJ2").Value = 7835167
$ws.Range("D3").Value = 2721
$ws.Range("E3").Value = 356
$ws.Range("F3").Value = 356
$ws.Range("G3").Value = 415
$ws.Range("H3").Value = 360
$ws.Range("I3").Value = 356
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 6985
$ws.Range("L3").Value = 1118
$ws.Range("M3").Value = 5868
$ws.Range("N3").Value = 5760
$ws.Range("O3").Value = 108
$ws.Range("P3").Value = 452
$ws.Range("Q3").Value = 71
$ws.Range("R3").Value = -152
$ws.Range("S3").Value = -29
$ws.Range("T3").Value = 138
$ws.Range("U3").Value = -68
$ws.Range("V3").Value = 396
$ws.Range("W3").Value = 13.08
$ws.Range("X3").Value = 13.25
$ws.Range("Y3").Value = 6.34
$ws.Range("Z3").Value = 5.24
$ws.Range("AA3").Value = 19.05
$ws.Range("AB3").Value = 1290.81
$ws.Range("AC3").Value = 4544
$ws.Range("AD3").Value = 9.13
$ws.Range("AE3").Value = 88729
$ws.Range("AF3").Value = 0.47
$ws.Range("AG3").Value = 1200
$ws.Range("AH3").Value = 2.89
$ws.Range("AI3").Value = 21.88
$ws.Range("AJ3").Value = 7835167
$ws.Range("D4").Value = 2976
$ws.Range("E4").Value = 665
$ws.Range("F4").Value = 665
$ws.Range("G4").Value = 684
$ws.Range("H4").Value = 612
$ws.Range("I4").Value = 609
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 7627
$ws.Range("L4").Value = 1193
$ws.Range("M4").Value = 6434
$ws.Range("N4").Value = 6333
$ws.Range("O4").Value = 101
$ws.Range("P4").Value = 452
$ws.Range("Q4").Value = 204
$ws.Range("R4").Value = -39
$ws.Range("S4").Value = -3
$ws.Range("T4").Value = 76
$ws.Range("U4").Value = 128
$ws.Range("V4").Value = 394
$ws.Range("W4").Value = 22.35
$ws.Range("X4").Value = 20.55
$ws.Range("Y4").Value = 10.07
$ws.Range("Z4").Value = 8.37
$ws.Range("AA4").Value = 18.54
$ws.Range("AB4").Value = 1412.8
$ws.Range("AC4").Value = 7769
$ws.Range("AD4").Value = 6.24
$ws.Range("AE4").Value = 97549
$ws.Range("AF4").Value = 0.5
$ws.Range("AG4").Value = 1400
$ws.Range("AH4").Value = 2.89
$ws.Range("AI4").Value = 14.93
$ws.Range("AJ4").Value = 7835167
$ws.Range("D5").Value = 3148
$ws.Range("E5").Value = 736
$ws.Range("F5").Value = 736
$ws.Range("G5").Value = 763
$ws.Range("H5").Value = 677
$ws.Range("I5").Value = 670
$ws.Range("J5").Value = 7
$ws.Range("K5").Value = 8185
$ws.Range("L5").Value = 1353
$ws.Range("M5").Value = 6832
$ws.Range("N5").Value = 6724
$ws.Range("O5").Value = 107
$ws.Range("P5").Value = 452
$ws.Range("Q5").Value = 220
$ws.Range("R5").Value = -287
$ws.Range("S5").Value = 41
$ws.Range("T5").Value = 187
$ws.Range("U5").Value = 34
$ws.Range("V5").Value = 430
$ws.Range("W5").Value = 23.39
$ws.Range("X5").Value = 21.51
$ws.Range("Y5").Value = 10.27
$ws.Range("Z5").Value = 8.56
$ws.Range("AA5").Value = 19.8
$ws.Range("AB5").Value = 1517.84
$ws.Range("AC5").Value = 8555
$ws.Range("AD5").Value = 6.14
$ws.Range("AE5").Value = 103582
$ws.Range("AF5").Value = 0.51
$ws.Range("AG5").Value = 1800
$ws.Range("AH5").Value = 3.43
$ws.Range("AI5").Value = 17.43
$ws.Range("AJ5").Value = 7835167
$ws.Range("D6").Value = 3061
$ws.Range("E6").Value = 365
$ws.Range("F6").Value = 365
$ws.Range("G6").Value = 405
$ws.Range("H6").Value = 355
$ws.Range("I6").Value = 350
$ws.Range("K6").Value = 8793
$ws.Range("L6").Value = 1659
$ws.Range("M6").Value = 7134
$ws.Range("N6").Value = 7033
$ws.Range("P6").Value = 452
$ws.Range("Q6").Value = -110
$ws.Range("R6").Value = -226
$ws.Range("S6").Value = 237
$ws.Range("T6").Value = 241
$ws.Range("U6").Value = -351
$ws.Range("V6").Value = 676
$ws.Range("W6").Value = 11.94
$ws.Range("X6").Value = 11.59
$ws.Range("Y6").Value = 5.09
$ws.Range("Z6").Value = 4.18
$ws.Range("AA6").Value = 23.26
$ws.Range("AB6").Value = 1584.19
$ws.Range("AC6").Value = 4469
$ws.Range("AD6").Value = 8.22
$ws.Range("AE6").Value = 108340
$ws.Range("AF6").Value = 0.34
$ws.Range("AG6").Value = 1700
$ws.Range("AH6").Value = 4.63
$ws.Range("AI6").Value = 31.52
$ws.Range("AJ6").Value = 7835167

# Clear the forecast rows (2019E, 2020E, 2021E) data cells, keep only A/B/C
$ws.Range("D7:AJ9").ClearContents()
